$d = $word.ActiveDocument

# Locate the paragraph containing the LOM3096 requirement line; the three
# paragraphs that follow it (a blank paragraph, "Ver no Jupiter..." line,
# and the "(c) 2020 ..." footer line) are being removed by this edit.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "LOM3096:") {
        $target = $i
        break
    }
}

if ($target -eq $null) {
    throw "Could not find LOM3096 paragraph"
}

$pBlank = $d.Paragraphs.Item($target + 1)
$pFooter2 = $d.Paragraphs.Item($target + 3)

$rng = $d.Range($pBlank.Range.Start, $pFooter2.Range.End)
$rng.Delete()
